# New crime data collected - weekly CompStat refresh (111th Precinct)
# Updates:
#   - Report "Number" (week-of-year counter) 28 -> 29
#   - Report date range 7/7/2025-7/13/2025 -> 7/14/2025-7/20/2025
#   - Crime-statistics table (rows 15-21, 24-27) refreshed with new weekly figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the "Volume 32   Number  28" banner text -> "...Number  29"
#    and the "Report Covering the Week 7/7/2025 Through 7/13/2025" banner
#    by editing only the affected substrings so the rest of the rich text
#    (fonts/sizes) is left untouched.
# ---------------------------------------------------------------------------

$numberCell = $ws.Range("A8")
$numberText = $numberCell.Text
$oldNum = "28"
$newNum = "29"
$idx = $numberText.IndexOf($oldNum)
$numberCell.Characters($idx + 1, $oldNum.Length).Text = $newNum

$weekCell = $ws.Range("C9")
$weekText = $weekCell.Text
$oldStart = "7/7/2025"
$newStart = "7/14/2025"
$startIdx = $weekText.IndexOf($oldStart)
$weekCell.Characters($startIdx + 1, $oldStart.Length).Text = $newStart

# re-read text since the string length changed after the first replacement
$weekText = $weekCell.Text
$oldEnd = "7/13/2025"
$newEnd = "7/20/2025"
$endIdx = $weekText.IndexOf($oldEnd)
$weekCell.Characters($endIdx + 1, $oldEnd.Length).Text = $newEnd

# ---------------------------------------------------------------------------
# 2. Helper to write a value into a given column of a given row.  Values
#    that start with "TXT:" are forced to be literal text (used for the "0"
#    and "***.*" placeholders that appear when a comparison period had zero
#    complaints, making a percentage change undefined).
# ---------------------------------------------------------------------------

function Set-StatCell($sheet, $row, $col, $value) {
    $cell = $sheet.Cells.Item($row, $col)
    if ($value -is [string] -and $value.StartsWith("TXT:")) {
        $text = $value.Substring(4)
        $cell.NumberFormat = "@"
        $cell.Value = $text
    } else {
        $cell.Value = $value
    }
}

$colIndex = @{ "C" = 3; "D" = 4; "E" = 5; "F" = 6; "G" = 7; "H" = 8; "I" = 9; "J" = 10; "K" = 11; "L" = 12; "M" = 13; "N" = 14 }

# ---------------------------------------------------------------------------
# 3. New weekly figures for each crime category row.
# ---------------------------------------------------------------------------

$rows = @(
    @{ Row = 15; C = 1;       D = "TXT:0";   E = "TXT:***.*"; F = 1;  G = 2;  H = -50;             I = 4;   J = 5;   K = -20;             L = -42.857142857142; M = 0;                 N = -42.857142857142 }
    @{ Row = 16; C = 1;       D = "TXT:0";   E = "TXT:***.*"; F = 2;  G = 5;  H = -60;             I = 16;  J = 33;  K = -51.515151515151; L = -57.894736842105; M = -65.957446808510; N = -91.208791208791 }
    @{ Row = 17; C = 1;       D = 2;         E = -50;         F = 11; G = 5;  H = 120;             I = 70;  J = 55;  K = 27.272727272727;  L = 27.272727272727;  M = 133.333333333333; N = 2.941176470588 }
    @{ Row = 18; C = 2;       D = 5;         E = -60;         F = 17; G = 17; H = 0;               I = 157; J = 141; K = 11.347517730496;  L = -8.720930232558;  M = 24.603174603174;  N = -72.310405643739 }
    @{ Row = 19; C = 11;      D = 12;        E = -8.333333333333; F = 35; G = 36; H = -2.777777777777; I = 244; J = 262; K = -6.870229007633;  L = -33.695652173913; M = 28.421052631578;  N = -15.277777777777 }
    @{ Row = 20; C = 8;       D = 5;         E = 60;          F = 23; G = 19; H = 21.052631578947;  I = 146; J = 140; K = 4.285714285714;   L = 69.767441860465;  M = 89.610389610389;  N = -91.732729331823 }
    @{ Row = 21; C = 24;      D = 24;        E = 0;           F = 89; G = 84; H = 5.952380952380;   I = 639; J = 636; K = 0.471698113207;   L = -12.225274725274; M = 34.810126582278;  N = -77.8125 }
    @{ Row = 24; C = 19;      D = 16;        E = 18.75;       F = 47; G = 51; H = -7.843137254901;  I = 307; J = 314; K = -2.229299363057;  L = -4.0625;           M = 18.076923076923;  N = "TXT:***.*" }
    @{ Row = 25; C = "TXT:0"; D = 5;         E = -100;        F = 4;  G = 12; H = -66.666666666666; I = 53;  J = 65;  K = -18.461538461538; L = -26.388888888888; M = "TXT:***.*";       N = "TXT:***.*" }
    @{ Row = 26; C = 6;       D = 6;         E = 0;           F = 24; G = 20; H = 20;               I = 115; J = 117; K = -1.709401709401;  L = 0.877192982456;   M = 29.213483146067;  N = "TXT:***.*" }
    @{ Row = 27; C = 1;       D = "TXT:0";   E = "TXT:***.*"; F = 1;  G = 3;  H = -66.666666666666; I = 6;   J = 8;   K = -25;              L = -25;              M = "TXT:***.*";       N = "TXT:***.*" }
)

foreach ($rowData in $rows) {
    $r = $rowData.Row
    foreach ($col in @("C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N")) {
        Set-StatCell $ws $r $colIndex[$col] $rowData[$col]
    }
}

Write-Host "CompStat weekly figures updated."
